$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: NODE -> NODE_ID (migrating remote_access import/export to v3)
$ws.Range("D1").Value = "NODE_ID"

# D column values switch from the text "master" to the numeric node id 1
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1

# Drop the one-off font override on A2 so it matches the sheet's normal font
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.Size = 11
$ws.Range("A2").Font.Color = 0

# Move the active selection to F3
$ws.Range("F3").Select() | Out-Null
